# Apply scheduled-runner updates to computed price/profit columns across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 279.73334
$ws.Range("J33").Value = 70
$ws.Range("L33").Value = 70
$ws.Range("N33").Value = -528

$ws.Range("H92").Value = 2955.2
$ws.Range("I92").Value = 1124.2142
$ws.Range("J92").Value = 5285.5454
$ws.Range("K92").Value = 1124.2142
$ws.Range("L92").Value = 5285.5454
$ws.Range("M92").Value = 123.7858000000001
$ws.Range("N92").Value = -7781.5454

$ws.Range("H99").Value = 3954.125
$ws.Range("I99").Value = 193
$ws.Range("K99").Value = 579
$ws.Range("M99").Value = 919

$ws.Range("H100").Value = 3508.1428
$ws.Range("I100").Value = 1771.2142
$ws.Range("J100").Value = 6982
$ws.Range("K100").Value = 1771.2142
$ws.Range("L100").Value = 6982
$ws.Range("M100").Value = -1230.2142
$ws.Range("N100").Value = -8064

$ws.Range("H104").Value = 166.75
$ws.Range("I104").Value = 93.5
$ws.Range("K104").Value = 280.5
$ws.Range("M104").Value = 1466.5

$ws.Range("H107").Value = 940.05884
$ws.Range("J107").Value = 2424.4
$ws.Range("L107").Value = 2424.4
$ws.Range("N107").Value = -6264.4

$ws.Range("H117").Value = 99999
$ws.Range("J117").Value = 99999
$ws.Range("L117").Value = 99999
$ws.Range("N117").Value = -109177

$ws.Range("H141").Value = 35722610
$ws.Range("I141").Value = 41670890
$ws.Range("K141").Value = 125012670
$ws.Range("M141").Value = -125007490

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 78178.164
$ws.Range("I60").Value = 83011.73
$ws.Range("K60").Value = 83011.73
$ws.Range("M60").Value = -82278.73

$ws.Range("H61").Value = 3091690
$ws.Range("I61").Value = 3281859.5
$ws.Range("J61").Value = 1434499.8
$ws.Range("K61").Value = 3281859.5
$ws.Range("L61").Value = 1434499.8
$ws.Range("M61").Value = -3281647.5
$ws.Range("N61").Value = -1434923.8

$ws.Range("H74").Value = 3125.8696
$ws.Range("I74").Value = 2471.3333
$ws.Range("K74").Value = 2471.3333
$ws.Range("M74").Value = -1597.3333

$ws.Range("H77").Value = 3125.8696
$ws.Range("I77").Value = 2471.3333
$ws.Range("K77").Value = 12356.6665
$ws.Range("M77").Value = -7988.666499999999

$ws.Range("H122").Value = 4075.8
$ws.Range("I122").Value = 3523.2727
$ws.Range("K122").Value = 10569.8181
$ws.Range("M122").Value = -8119.8181

$ws.Range("H136").Value = 3091690
$ws.Range("I136").Value = 3281859.5
$ws.Range("J136").Value = 1434499.8
$ws.Range("K136").Value = 9845578.5
$ws.Range("L136").Value = 4303499.4
$ws.Range("M136").Value = -9843028.5
$ws.Range("N136").Value = -4308599.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11766.889
$ws.Range("I20").Value = 18559
$ws.Range("J20").Value = 3276.75
$ws.Range("K20").Value = 18559
$ws.Range("L20").Value = 3276.75
$ws.Range("M20").Value = -18312
$ws.Range("N20").Value = -3770.75

$ws.Range("H107").Value = 2779
$ws.Range("I107").Value = 2944.9092
$ws.Range("K107").Value = 2944.9092
$ws.Range("M107").Value = -1024.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1353.3572
$ws.Range("J94").Value = 1426.3334
$ws.Range("L94").Value = 1426.3334
$ws.Range("N94").Value = -2328.3334

$ws.Range("H105").Value = 1415
$ws.Range("I105").Value = 947.5
$ws.Range("J105").Value = 1602
$ws.Range("K105").Value = 947.5
$ws.Range("L105").Value = 1602
$ws.Range("M105").Value = 799.5
$ws.Range("N105").Value = -5096

$ws.Range("H132").Value = 1402.6428
$ws.Range("I132").Value = 1402.6428
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4207.928400000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1677.928400000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 9391.5
$ws.Range("I55").Value = 2160.2
$ws.Range("K55").Value = 6480.599999999999
$ws.Range("M55").Value = -6303.599999999999

$ws.Range("H57").Value = 13822.667
$ws.Range("I57").Value = 7505
$ws.Range("J57").Value = 14794.615
$ws.Range("K57").Value = 22515
$ws.Range("L57").Value = 44383.845
$ws.Range("M57").Value = -21956
$ws.Range("N57").Value = -45501.845

$ws.Range("H59").Value = 9683.429
$ws.Range("I59").Value = 4650.6665
$ws.Range("J59").Value = 13458
$ws.Range("K59").Value = 13951.9995
$ws.Range("L59").Value = 40374
$ws.Range("M59").Value = -13411.9995
$ws.Range("N59").Value = -41454

$ws.Range("H60").Value = 4522.9
$ws.Range("I60").Value = 275
$ws.Range("J60").Value = 7354.8335
$ws.Range("K60").Value = 825
$ws.Range("L60").Value = 22064.5005
$ws.Range("M60").Value = -574
$ws.Range("N60").Value = -22566.5005

$ws.Range("H63").Value = 20776.572
$ws.Range("J63").Value = 26032.7
$ws.Range("L63").Value = 78098.10000000001
$ws.Range("N63").Value = -79596.10000000001

$ws.Range("H64").Value = 13585.714
$ws.Range("I64").Value = 5687
$ws.Range("J64").Value = 33332.5
$ws.Range("K64").Value = 17061
$ws.Range("L64").Value = 99997.5
$ws.Range("M64").Value = -16791
$ws.Range("N64").Value = -100537.5

$ws.Range("H66").Value = 20776.572
$ws.Range("J66").Value = 26032.7
$ws.Range("L66").Value = 234294.3
$ws.Range("N66").Value = -241782.3

$ws.Range("H67").Value = 13585.714
$ws.Range("I67").Value = 5687
$ws.Range("J67").Value = 33332.5
$ws.Range("K67").Value = 17061
$ws.Range("L67").Value = 99997.5
$ws.Range("M67").Value = -16125
$ws.Range("N67").Value = -101869.5

$ws.Range("H68").Value = 2774.5
$ws.Range("I68").Value = 550
$ws.Range("K68").Value = 1650
$ws.Range("M68").Value = -839

$ws.Range("H70").Value = 17223.166
$ws.Range("I70").Value = 12501.5
$ws.Range("J70").Value = 26666.5
$ws.Range("K70").Value = 37504.5
$ws.Range("L70").Value = 79999.5
$ws.Range("M70").Value = -37189.5
$ws.Range("N70").Value = -80629.5

$ws.Range("H71").Value = 2774.5
$ws.Range("I71").Value = 550
$ws.Range("K71").Value = 4950
$ws.Range("M71").Value = -894

$ws.Range("H73").Value = 17223.166
$ws.Range("I73").Value = 12501.5
$ws.Range("J73").Value = 26666.5
$ws.Range("K73").Value = 37504.5
$ws.Range("L73").Value = 79999.5
$ws.Range("M73").Value = -36412.5
$ws.Range("N73").Value = -82183.5

$ws.Range("H82").Value = 17500
$ws.Range("I82").Value = 17500
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 52500
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -52094
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 17500
$ws.Range("I85").Value = 17500
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 52500
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -51096
$ws.Range("N85").ClearContents()

$ws.Range("H110").Value = 18775
$ws.Range("I110").Value = 4217
$ws.Range("K110").Value = 12651
$ws.Range("M110").Value = -8561

$ws.Range("H131").Value = 4781.1904
$ws.Range("I131").Value = 3430
$ws.Range("J131").Value = 5612.6924
$ws.Range("K131").Value = 10290
$ws.Range("L131").Value = 16838.0772
$ws.Range("M131").Value = -5250
$ws.Range("N131").Value = -26918.0772

$ws.Range("H132").Value = 1907.6842
$ws.Range("I132").Value = 1803.5
$ws.Range("K132").Value = 16231.5
$ws.Range("M132").Value = -13701.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5035.4814
$ws.Range("I122").Value = 4513.5713
$ws.Range("J122").Value = 6862.1665
$ws.Range("K122").Value = 13540.7139
$ws.Range("L122").Value = 20586.4995
$ws.Range("M122").Value = -11090.7139
$ws.Range("N122").Value = -25486.4995

$ws.Range("H131").Value = 89176.8
$ws.Range("J131").Value = 89176.8
$ws.Range("L131").Value = 89176.8
$ws.Range("N131").Value = -99256.8

$ws.Range("H132").Value = 3573129.2
$ws.Range("I132").Value = 1592.32
$ws.Range("K132").Value = 4776.96
$ws.Range("M132").Value = -2246.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8690.714
$ws.Range("I22").Value = 16818.428
$ws.Range("K22").Value = 16818.428
$ws.Range("M22").Value = -16523.428

$ws.Range("H27").Value = 8690.714
$ws.Range("I27").Value = 16818.428
$ws.Range("K27").Value = 16818.428
$ws.Range("M27").Value = -16711.428

$ws.Range("H46").Value = 1068.2
$ws.Range("I46").Value = 830.3333
$ws.Range("K46").Value = 830.3333
$ws.Range("M46").Value = -642.3333

$ws.Range("H61").Value = 3745.4644
$ws.Range("I61").Value = 3057
$ws.Range("J61").Value = 5810.857
$ws.Range("K61").Value = 3057
$ws.Range("L61").Value = 5810.857
$ws.Range("M61").Value = -2855
$ws.Range("N61").Value = -6214.857

$ws.Range("H68").Value = 3208197.2
$ws.Range("I68").Value = 5954651
$ws.Range("J68").Value = 4001.3333
$ws.Range("K68").Value = 5954651
$ws.Range("L68").Value = 4001.3333
$ws.Range("M68").Value = -5953902
$ws.Range("N68").Value = -5499.3333

$ws.Range("H71").Value = 3208197.2
$ws.Range("I71").Value = 5954651
$ws.Range("J71").Value = 4001.3333
$ws.Range("K71").Value = 29773255
$ws.Range("L71").Value = 20006.6665
$ws.Range("M71").Value = -29769511
$ws.Range("N71").Value = -27494.6665

$ws.Range("H113").Value = 3745.4644
$ws.Range("I113").Value = 3057
$ws.Range("J113").Value = 5810.857
$ws.Range("K113").Value = 3057
$ws.Range("L113").Value = 5810.857
$ws.Range("M113").Value = -887
$ws.Range("N113").Value = -10150.857

$ws.Range("H132").Value = 3977.9524
$ws.Range("I132").Value = 2545.2222
$ws.Range("J132").Value = 5052.5
$ws.Range("K132").Value = 7635.6666
$ws.Range("L132").Value = 15157.5
$ws.Range("M132").Value = -5105.6666
$ws.Range("N132").Value = -20217.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 13750
$ws.Range("I31").Value = 13750
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 13750
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -13402
$ws.Range("N31").ClearContents()
